# Update GUI Plotting for Multiple Scenarios
# Renames several y_data label strings on the plot_definition sheet so that
# the "hauling/non hauling" + "BEV/ICE" qualifier is expressed as a suffix
# on the metric name (average_cost_... / average_co2_gpmi_...) instead of a
# prefix, and adjusts column widths / selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Average Vehicle CO2 plot): co2_gpmi columns ---
$ws.Range("I2").Value = "average_co2_gpmi_non hauling.BEV"
$ws.Range("J2").Value = "average_co2_gpmi_hauling.BEV"
$ws.Range("K2").Value = "average_co2_gpmi_non hauling.ICE"
$ws.Range("L2").Value = "average_co2_gpmi_hauling.ICE"

# --- Row 3 (Average Vehicle Cost plot): cost columns ---
$ws.Range("G3").Value = "average_cost_non hauling.ICE"
$ws.Range("I3").Value = "average_cost_non hauling.BEV"
$ws.Range("J3").Value = "average_cost_hauling.BEV"
$ws.Range("K3").Value = "average_cost_non hauling.BEV"
$ws.Range("L3").Value = "average_cost_hauling.ICE"

# --- Column widths (engine rounds ColumnWidth to a 1/6-character pixel grid,
#     so the values below are chosen to land as close as possible to the
#     widths recorded in the saved file) ---
$ws.Columns.Item(7).ColumnWidth = 26.5
$ws.Columns.Item(9).ColumnWidth = 38
$ws.Columns.Item(11).ColumnWidth = 31.666666666666668
$ws.Columns.Item(12).ColumnWidth = 25.833333333333332

# --- Selection moved from the data rows to the (now blank) rows below ---
$ws.Range("A9:XFD12").Select()
